# Lab 1/Lab 2 block diagram - final minor corrections:
#   1. Bump the cached "datetimeFigureOut" date placeholder text from
#      2/6/2017 to 2/8/2017 everywhere it is cached (slide master, every
#      slide layout, and the notes master).
#   2. Remove the stray duplicate "clk" label text box on the slide.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq "2/6/2017") {
                $sh.TextFrame.TextRange.Text = "2/8/2017"
            }
        }
    }
}

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes

# Every slide layout off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Set-DatePlaceholderText $layouts.Item($L).Shapes
}

# Notes master
Set-DatePlaceholderText $p.NotesMaster.Shapes

# Remove the stray "clk" text box (TextBox 429) from slide 1
$s = $p.Slides.Item(1)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 429") {
        $sh.Delete()
    }
}
